# Swap the presentation's active theme (ppt/theme/theme1.xml, the theme
# referenced by the slide master / all slides) from the "Integral" theme
# palette to the stock "Office Theme" palette.
#
# Both themes share an identical font scheme ("Office": Arial everywhere)
# and an identical format scheme (fill/line/effect/background styles), so
# the only substantive difference between them is the 12-colour theme
# colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink). We rewrite
# those twelve colours in place via ThemeColorScheme, which is the COM
# surface that maps onto <a:clrScheme> in ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette, in ThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $officeHex.Length; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB() values are packed 0x00BBGGRR
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $tcs.Colors($i).RGB = $bgr
}
